$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Epayco sheet: add new row 3 (F3 = 1077)
# ---------------------------------------------------------------------------
$epayco = $wb.Worksheets.Item("Epayco")
$epayco.Range("F3").Value = 1077

# ---------------------------------------------------------------------------
# 2) Equipo sheet: widen column A and move the selection to D6
# ---------------------------------------------------------------------------
$equipo = $wb.Worksheets.Item("Equipo")
$equipo.Columns.Item(1).ColumnWidth = 11.666666666666666
$equipo.Activate() | Out-Null
$equipo.Range("D6").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3) Add the new "Ambiente" sheet at the end of the workbook
# ---------------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$ambiente = $wb.Worksheets.Add($null, $lastSheet)
$ambiente.Name = "Ambiente"

# Write the cell values in the same order the original authoring tool used
# so new shared-strings get appended in a matching order.
$ambiente.Range("A1").Value = "URL Pruebas"
$ambiente.Range("B1").Value = "URL Produccion"
$ambiente.Range("A2").Value = "https://integration-5ojmyuq-jvrr247te2phq.us-5.magentosite.cloud/celulares.html"
$ambiente.Range("B2").Value = "https://tienda.movistar.com.co/celulares.html"
$ambiente.Range("C1").Value = "URL"
$ambiente.Range("C2").Value = "https://integration-5ojmyuq-jvrr247te2phq.us-5.magentosite.cloud/celulares.html"

# Hyperlinks on A2/B2 (this also gives those two cells the hyperlink style)
$ambiente.Hyperlinks.Add($ambiente.Range("A2"), "https://integration-5ojmyuq-jvrr247te2phq.us-5.magentosite.cloud/celulares.html") | Out-Null
$ambiente.Hyperlinks.Add($ambiente.Range("B2"), "https://tienda.movistar.com.co/celulares.html") | Out-Null

# Drop-down list validation on C2, sourced from A2:B2
$ambiente.Range("C2").Validation.Add(3, 1, 1, "=`$A`$2:`$B`$2") | Out-Null

# Column widths
$ambiente.Columns.Item(1).ColumnWidth = 15.0
$ambiente.Columns.Item(2).ColumnWidth = 15.666666666666666
$ambiente.Columns.Item(3).ColumnWidth = 68.0

# Selection left on the sheet after editing
$ambiente.Range("C4").Select() | Out-Null
